{"js": "const pairs = [\n  [\"854\u00f76=\", \"786\u00f76=\"],\n  [\"226\u00f75=\", \"843\u00f76=\"],\n  [\"863\u00f73=\", \"617\u00f72=\"],\n  [\"972\u00f79=\", \"566\u00f76=\"],\n  [\"162\u00f79=\", \"449\u00f78=\"],\n  [\"651\u00f73=\", \"463\u00f73=\"],\n  [\"501\u00f72=\", \"970\u00f73=\"],\n  [\"352\u00f76=\", \"861\u00f79=\"],\n  [\"953\u00f74=\", \"417\u00f73=\"],\n  [\"745\u00f79=\", \"622\u00f78=\"],\n  [\"896\u00f79=\", \"935\u00f77=\"],\n  [\"995\u00f77=\", \"178\u00f75=\"],\n  [\"222\u00f78=\", \"939\u00f76=\"],\n  [\"915\u00f74=\", \"239\u00f73=\"],\n  [\"148\u00f79=\", \"739\u00f77=\"],\n  [\"706\u00f76=\", \"172\u00f77=\"],\n  [\"190\u00f74=\", \"161\u00f72=\"],\n  [\"153\u00f77=\", \"421\u00f79=\"],\n  [\"139\u00f72=\", \"519\u00f72=\"],\n  [\"724\u00f73=\", \"787\u00f78=\"],\n  [\"196\u00f74=\", \"146\u00f72=\"],\n  [\"791\u00f77=\", \"679\u00f72=\"],\n  [\"710\u00f75=\", \"998\u00f76=\"],\n  [\"213\u00f79=\", \"488\u00f78=\"],\n  [\"271\u00f76=\", \"869\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"854\u00f76=\", \"786\u00f76=\"),\n    @(\"226\u00f75=\", \"843\u00f76=\"),\n    @(\"863\u00f73=\", \"617\u00f72=\"),\n    @(\"972\u00f79=\", \"566\u00f76=\"),\n    @(\"162\u00f79=\", \"449\u00f78=\"),\n    @(\"651\u00f73=\", \"463\u00f73=\"),\n    @(\"501\u00f72=\", \"970\u00f73=\"),\n    @(\"352\u00f76=\", \"861\u00f79=\"),\n    @(\"953\u00f74=\", \"417\u00f73=\"),\n    @(\"745\u00f79=\", \"622\u00f78=\"),\n    @(\"896\u00f79=\", \"935\u00f77=\"),\n    @(\"995\u00f77=\", \"178\u00f75=\"),\n    @(\"222\u00f78=\", \"939\u00f76=\"),\n    @(\"915\u00f74=\", \"239\u00f73=\"),\n    @(\"148\u00f79=\", \"739\u00f77=\"),\n    @(\"706\u00f76=\", \"172\u00f77=\"),\n    @(\"190\u00f74=\", \"161\u00f72=\"),\n    @(\"153\u00f77=\", \"421\u00f79=\"),\n    @(\"139\u00f72=\", \"519\u00f72=\"),\n    @(\"724\u00f73=\", \"787\u00f78=\"),\n    @(\"196\u00f74=\", \"146\u00f72=\"),\n    @(\"791\u00f77=\", \"679\u00f72=\"),\n    @(\"710\u00f75=\", \"998\u00f76=\"),\n    @(\"213\u00f79=\", \"488\u00f78=\"),\n    @(\"271\u00f76=\", \"869\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $pair[0]\n    $range.Find.Replacement.Text = $pair[1]\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
